$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- User profile block (rows 5-9), columns: C=User1, E=User2, F=User3, G=User4 ---
# Order follows the original authoring sequence: User1, User4, User2, User3

# User1 (column C)
$ws.Range("C5").Value = "Hombre"
$ws.Range("C7").Value = "Funcionario"
$ws.Range("C8").Value = "Intermedio"
$ws.Range("C9").Value = "El usuario 1 es un padre de familia de mediana edad que trabaja de funcionario. Le gusta mucho estudiar idiomas y se encuentra algo triste"

# User4 (column G)
$ws.Range("G5").Value = "Mujer"
$ws.Range("G7").Value = "Gerente"
$ws.Range("G8").Value = "Bajo"
$ws.Range("G9").Value = "El usuario 4 es una mujer descapacitada bastante planificadora y tiene algo de miedo"

# User2 (column E)
$ws.Range("E5").Value = "Mujer "
$ws.Range("E7").Value = "Estudiante"
$ws.Range("E8").Value = "Avanzado"
$ws.Range("E9").Value = "Estudiante con una discapacidad de movilidad que se encuentra disgustada"

# User3 (column F)
$ws.Range("F5").Value = "Hombre"
$ws.Range("F7").Value = "Becario"
$ws.Range("F8").Value = "Avanzado"
$ws.Range("F9").Value = "El usuario 3 es una persona sociable, interesado en aprender idiomas y está sorprendido"

# Ages (row 6) -- numeric
$ws.Range("C6").Value = 46
$ws.Range("E6").Value = 23
$ws.Range("F6").Value = 28
$ws.Range("G6").Value = 60

# --- SUS questionnaire answers (rows 14-23), columns C/E/F/G ---
$ws.Range("C14").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1

$ws.Range("C15").Value = 2
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3

$ws.Range("C16").Value = 5
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 3

$ws.Range("C17").Value = 1
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1

$ws.Range("C18").Value = 4
$ws.Range("E18").Value = 5
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 4

$ws.Range("C19").Value = 2
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 1

$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 4
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 4

$ws.Range("C21").Value = 4
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 3

$ws.Range("C22").Value = 4
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 3

$ws.Range("C23").Value = 1
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 1
